# mobiletestinglogins.xlsx - "updated with datadriven approch for testing multiple set of data"
#
# The real content edit is two corrected email addresses on the DATA_SHEET
# tab (column C, the EMAIL column):
#   C2: DUVARAKESH123@GMAIL.COM -> DUVAKSH123@GMAIL.COM
#   C5: DUVAR3212@GMAIL.COM     -> AR321882@GMAIL.COM

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("DATA_SHEET")
$outputSheet = $wb.Worksheets.Item("OUTPUT_DATASHEET")

$dataSheet.Range("C2").Value = "DUVAKSH123@GMAIL.COM"
$dataSheet.Range("C5").Value = "AR321882@GMAIL.COM"

# Match the author's final selection/active-sheet state: DATA_SHEET ends up
# the active tab with C5 selected, while OUTPUT_DATASHEET's last selection
# was A10.
$outputSheet.Range("A10").Select()
$dataSheet.Activate()
$dataSheet.Range("C5").Select()
